# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two label-only "section header" rows (row 5 "situação do
# domicílio" and row 8 "grandes regiões e unidades da federação") that
# carried no data — an artifact that also left every data row below them
# misaligned by one slot relative to its real label. Removing the two
# empty rows (and letting the rows below shift up) fixes the label/data
# alignment for "urbana"/"rural" and all the region/state rows, and
# shrinks the used range from A1:F39 to A1:F37.
#
# Row 2's sub-header cells (B2/F2) were placeholder "unnamed: ..." labels
# that should simply read "total", matching C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two empty separator rows. Row 8 is deleted first (higher
# index) so row 5's index is still valid when it's deleted afterwards.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Fix the row-2 sub-header labels.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
